$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.852.21'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").Value = '1.812.28'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5908'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2753'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06740'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.30%  '
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("D12").Value = '1.818.00'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.669'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6239'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009311'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '74.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.00%  '
$ws.Range("D17").Value = '28.615.02'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.439'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.70%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.764'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.05%  '
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1271'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.791'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06297'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.30%  '
$ws.Range("E29").Value = '  -5.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.427'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.728'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.686'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.693'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("E34").Value = '  -7.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6347'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.52%  '
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.727'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.426'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("E39").Value = '  -4.36%  '
$ws.Range("D40").Value = '1.131.01'
$ws.Range("E40").Value = '  -8.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8674'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.15%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = '1.968.48'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.05%  '
$ws.Range("E46").Value = '  -3.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.569'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4512'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05454'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.263'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '
